# Generate Report for Handoff
# Updates the localization-status report to reflect that
# ba847e91-09e3-46c3-9e32-671142e51268.md is now "Ready for handoff",
# and refreshes the "Latest Handoff" timestamps for the rows whose
# handoff data changed (1a0d14b1... and ba847e91... -> cd733652... region).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("D6").Value = "2016-03-24 12:34:20"

$ov.Range("B9").Value = "Ready for handoff"
$ov.Range("C9").Value = "Ready for handoff"
$ov.Range("D9").Value = "2016-03-24 12:34:20"

$ov.Range("D10").Value = "2016-03-24 12:34:20"

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("E6").Value = "2016-03-24 12:34:15"

$zh.Range("C9").Value = "Ready for handoff"
$zh.Range("E9").Value = "2016-03-24 12:34:15"

$zh.Range("E10").Value = "2016-03-24 12:34:15"

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("E6").Value = "2016-03-24 12:34:20"

$de.Range("C9").Value = "Ready for handoff"
$de.Range("E9").Value = "2016-03-24 12:34:20"

$de.Range("E10").Value = "2016-03-24 12:34:20"
